# Updates the cryptocurrency price list (columns D = Price, E = Volume(1h))
# with freshly scraped values, as produced by the scheduled GitHub Actions
# "Updated cryptos list" job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-wise updates: each entry carries the new Price (column D) and the new
# Volume(1h) percentage text (column E). A $null member means that column is
# left untouched for that row.
$updates = @(
    @{ Row = 2; D = "28.027.68"; E = "  -0.66%  " },
    @{ Row = 3; D = "1.900.29"; E = "  +1.58%  " },
    @{ Row = 4; D = "1.002"; E = "  +0.08%  " },
    @{ Row = 5; D = "312.49"; E = "  +0.24%  " },
    @{ Row = 7; D = "0.5082"; E = "  +0.62%  " },
    @{ Row = 8; D = "0.3919"; E = "  -0.03%  " },
    @{ Row = 9; D = "0.09222"; E = "  -4.42%  " },
    @{ Row = 10; D = "1.133"; E = "  -0.52%  " },
    @{ Row = 11; D = "41.75"; E = "  +2.16%  " },
    @{ Row = 12; D = "6.356"; E = "  -2.18%  " },
    @{ Row = 13; D = "20.75"; E = "  -0.98%  " },
    @{ Row = 14; D = "1.898.04"; E = "  +1.19%  " },
    @{ Row = 15; D = "1.002"; E = "  +0.13%  " },
    @{ Row = 16; D = "7.285"; E = "  -1.89%  " },
    @{ Row = 17; D = "0.00001115"; E = "  -1.13%  " },
    @{ Row = 18; D = "92.31"; E = "  -0.67%  " },
    @{ Row = 19; D = "0.06577"; E = "  -0.81%  " },
    @{ Row = 20; D = "17.76"; E = "  +1.17%  " },
    @{ Row = 21; D = $null; E = "  -0.01%  " },
    @{ Row = 22; D = "6.214"; E = "  +0.94%  " },
    @{ Row = 23; D = "28.091.06"; E = "  -0.61%  " },
    @{ Row = 24; D = "11.34"; E = "  -0.04%  " },
    @{ Row = 25; D = "2.319"; E = "  +1.64%  " },
    @{ Row = 26; D = "2.602"; E = "  +2.71%  " },
    @{ Row = 27; D = "2.115.99"; E = "  +1.17%  " },
    @{ Row = 28; D = "20.89"; E = "  -1.47%  " },
    @{ Row = 29; D = "157.51"; E = "  +0.08%  " },
    @{ Row = 30; D = "127.12"; E = "  -0.15%  " },
    @{ Row = 31; D = "1.085"; E = "  +1.54%  " },
    @{ Row = 32; D = $null; E = "  +1.02%  " },
    @{ Row = 33; D = "5.600"; E = "  -0.42%  " },
    @{ Row = 34; D = "3.611"; E = "  -0.49%  " },
    @{ Row = 35; D = "9.588"; E = "  +0.16%  " },
    @{ Row = 36; D = "0.06641"; E = "  -1.55%  " },
    @{ Row = 37; D = "0.02403"; E = "  +0.73%  " },
    @{ Row = 38; D = $null; E = "  -0.51%  " },
    @{ Row = 39; D = "1.222"; E = "  -1.63%  " },
    @{ Row = 40; D = "1.257"; E = "  +6.68%  " },
    @{ Row = 41; D = "0.6356"; E = "  +0.10%  " },
    @{ Row = 42; D = "4.974"; E = "  -0.25%  " },
    @{ Row = 43; D = "11.39"; E = "  -0.69%  " },
    @{ Row = 44; D = $null; E = "  +0.02%  " },
    @{ Row = 45; D = "13.28"; E = $null },
    @{ Row = 46; D = "0.5975"; E = "  -0.73%  " },
    @{ Row = 47; D = "3.701"; E = "  +1.04%  " },
    @{ Row = 48; D = "1.276"; E = "  +0.83%  " },
    @{ Row = 49; D = "2.006"; E = "  +0.70%  " },
    @{ Row = 50; D = "122.40"; E = "  -1.46%  " },
    @{ Row = 51; D = "1.177"; E = "  -1.51%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Range("D$row")
        # Several prices look like plain numbers (or dates once Excel parses
        # the dots), so force the cell to Text format before assigning and
        # restore its original style afterwards to avoid altering formatting.
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = $origStyle
    }

    if ($null -ne $u.E) {
        # Volume text already carries surrounding spaces + a percent sign,
        # so Excel keeps it as plain text without any extra handling.
        $ws.Range("E$row").Value = $u.E
    }
}
